$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fees")

# Row 23: new timesheet entry (Factor Analysis for 2013 MTF future outlook, social media,
# parental background groups) added underneath the existing 12/7/2020 entry, copying the
# formatting of the row above it (row 22).
$ws.Rows("23").RowHeight = 30

$ws.Range("B23").Value = 44172

$ws.Range("C23").Value = "Factor Analysis for 2013 MTF future outlook, social media, parental background groups"
$ws.Range("C23").WrapText = $true

$ws.Range("D23").Value = 2.5
$ws.Range("D23").NumberFormat = "0.0"

$ws.Range("E23").Formula = "=D23*40"
$ws.Range("E23").NumberFormat = '_($* #,##0.00_);_($* (#,##0.00);_($* "-"??_);_(@_)'

# Move the active selection to match the author's final cursor position.
$ws.Range("G25").Select()

$wb.Save()
